$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are written as Text so that values
# like "1.00" or "8.20" are not silently coerced into numbers,
# matching the original inlineStr/text storage of this column.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.716.55"
$ws.Range("E2").Value = "  -6.03%  "

# Row 3
$ws.Range("D3").Value = "3.263.35"
$ws.Range("E3").Value = "  -8.21%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("D5").Value = "176.23"
$ws.Range("E5").Value = "  -15.98%  "

# Row 6
$ws.Range("D6").Value = "511.13"
$ws.Range("E6").Value = "  -9.46%  "

# Row 7
$ws.Range("D7").Value = "0.586"
$ws.Range("E7").Value = "  -4.02%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.255.59"
$ws.Range("E9").Value = "  -8.31%  "

# Row 10
$ws.Range("D10").Value = "0.610"
$ws.Range("E10").Value = "  -9.26%  "

# Row 11
$ws.Range("D11").Value = "56.63"
$ws.Range("E11").Value = "  -7.41%  "

# Row 12
$ws.Range("D12").Value = "0.129"
$ws.Range("E12").Value = "  -11.66%  "

# Row 13
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  -9.09%  "

# Row 14
$ws.Range("D14").Value = "8.98"
$ws.Range("E14").Value = "  -11.43%  "

# Row 15
$ws.Range("D15").Value = "3.806.42"
$ws.Range("E15").Value = "  -8.01%  "

# Row 16
$ws.Range("D16").Value = "0.117"
$ws.Range("E16").Value = "  -7.01%  "

# Row 17
$ws.Range("D17").Value = "3.284.78"
$ws.Range("E17").Value = "  -7.77%  "

# Row 18
$ws.Range("D18").Value = "63.573.31"
$ws.Range("E18").Value = "  -6.03%  "

# Row 19
$ws.Range("D19").Value = "16.99"
$ws.Range("E19").Value = "  -10.21%  "

# Row 20
$ws.Range("D20").Value = "10.69"
$ws.Range("E20").Value = "  -11.61%  "

# Row 21
$ws.Range("D21").Value = "0.934"
$ws.Range("E21").Value = "  -11.29%  "

# Row 22
$ws.Range("D22").Value = "365.99"
$ws.Range("E22").Value = "  -8.35%  "

# Row 23
$ws.Range("D23").Value = "79.24"
$ws.Range("E23").Value = "  -5.34%  "

# Row 24
$ws.Range("D24").Value = "3.60"
$ws.Range("E24").Value = "  -12.87%  "

# Row 25
$ws.Range("D25").Value = "10.76"
$ws.Range("E25").Value = "  -14.36%  "

# Row 26
$ws.Range("D26").Value = "3.75"
$ws.Range("E26").Value = "  -3.63%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.61"
$ws.Range("E27").Value = "  -8.47%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "11.09"
$ws.Range("E28").Value = "  -9.89%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "8.20"
$ws.Range("E29").Value = "  -10.24%  "

# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "645.93"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31
$ws.Range("D31").Value = "28.06"
$ws.Range("E31").Value = "  -10.24%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "6.57"
$ws.Range("E32").Value = "  -13.80%  "

# Row 33
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "10.94"
$ws.Range("E33").Value = "  -8.73%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "58.61"
$ws.Range("E34").Value = "  -6.81%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.102"
$ws.Range("E35").Value = "  -8.94%  "

# Row 36
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "35.45"
$ws.Range("E37").Value = "  -13.07%  "

# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "0.372"
$ws.Range("E38").Value = "  -8.37%  "

# Row 39
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.121"
$ws.Range("E40").Value = "  -8.30%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.835.52"
$ws.Range("E41").Value = "  -10.10%  "

# Row 42
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0640"
$ws.Range("E42").Value = "  -13.47%  "

# Row 43
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  -18.51%  "

# Row 44
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  -6.66%  "

# Row 45
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.32"
$ws.Range("E45").Value = "  -12.69%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0377"
$ws.Range("E46").Value = "  -7.58%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  +2.12%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "2.91"
$ws.Range("E48").Value = "  -3.04%  "

# Row 49
$ws.Range("D49").Value = "0.122"
$ws.Range("E49").Value = "  -5.48%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "131.74"
$ws.Range("E50").Value = "  -4.38%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "2.37"
$ws.Range("E51").Value = "  -6.15%  "
